# The "2010" sheet (sheet1.xml) had a blank/unused column B between the
# "Uraian" label column (A) and the "Kab. Badung" data column (formerly C).
# This removes that empty column, shifting all data (Kab. Badung ... Kota
# Denpasar, formerly columns C:K) one column to the left (B:J), which
# matches the fixed/re-imported layout from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("B").Delete()

# Restore the cursor/selection position left behind by the edit.
$ws.Range("E8").Select() | Out-Null
